$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 7 (pushes existing rows 7..67 down to 8..68)
$ws.Rows.Item(7).Insert()

# Populate the new row 7 with the new record's data
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C7").Value = "Metropolitana"
$ws.Range("D7").Value = 44630
$ws.Range("E7").Value = 13
$ws.Range("F7").Value = 100114007
$ws.Range("G7").Value = "Jengibre"
$ws.Range("H7").Value = "Sin especificar"
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 220
$ws.Range("K7").Value = 13000
$ws.Range("L7").Value = 15000
$ws.Range("M7").Value = 13909
$ws.Range("N7").Value = "`$/caja 13 kilos"
$ws.Range("O7").Value = "Perú"
$ws.Range("P7").Value = 1070
$ws.Range("Q7").Value = 13
$ws.Range("R7").Value = "Hortaliza"
